$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing 20 data rows (rows 2..21) down by 2, to rows 4..23.
# Walk from the bottom up so we never clobber a row before reading it.
for ($r = 21; $r -ge 2; $r--) {
    $dstRow = $r + 2
    $ws.Cells.Item($dstRow, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dstRow, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dstRow, 3).Value = $ws.Cells.Item($r, 3).Value2
}

# Write the two brand-new rows that now occupy rows 2-3
$ws.Cells.Item(2, 1).Value = 0.001320064067840854
$ws.Cells.Item(2, 2).Value = 0.1883212360553446
$ws.Cells.Item(2, 3).Value = -0.001869207888376141

$ws.Cells.Item(3, 1).Value = -0.1889566183090211
$ws.Cells.Item(3, 2).Value = 0.01437168661504978
$ws.Cells.Item(3, 3).Value = 0.113000919460319

# Append 8 new rows of data at the bottom (rows 24-31)
$bottom = @(
    @(4.677844420075353, -3.651133604347696, -7.842656075954431),
    @(-2.627654522657398, -2.928949266672134, 4.230176210403448),
    @(-4.852406792342663, 0.3913787733763447, 0.2968738228082666),
    @(-1.301035702228551, 3.64691380783915, -6.109266191720954),
    @(2.465943455696097, -2.991184197366218, -3.608212560415278),
    @(-1.307898223400096, -2.068972408771528, -0.7334359884262174),
    @(-1.702915767207749, -0.5735956337302961, -0.9715757742524092),
    @(-3.414293382316824, 0.2869436666369428, -0.1008520126342796)
)

for ($i = 0; $i -lt $bottom.Length; $i++) {
    $row = 24 + $i
    $ws.Cells.Item($row, 1).Value = $bottom[$i][0]
    $ws.Cells.Item($row, 2).Value = $bottom[$i][1]
    $ws.Cells.Item($row, 3).Value = $bottom[$i][2]
}
